# Deploying to gh-pages: add the 2020 data column (M) to the sanitation-access
# table and refresh the saved selection, mirroring the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant: copying the neighbouring cell's style onto the new
# M-column cell so the new data matches the look of the rest of its row
# (this also leaves the old "empty-placeholder" style unused on save, just
# like the source edit did).
$xlPasteFormats = -4122

# --- Row 3: extend the thin top border strip one cell to the right (M3) ---
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial($xlPasteFormats)

# --- Row 4: new "2020" column header, formatted like the other year cells ---
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial($xlPasteFormats)
$ws.Range("M4").Value = 2020

# --- Rows 5-14: the 2020 data points for each region/total line ---
$dataRows = @{
    5  = 34.377950588852634
    6  = 4.8358243107925931
    7  = 5.9543034993102522
    8  = 51.21106605430419
    9  = 27.156801192263725
    10 = 0.94331159862228353
    11 = 7.8509592890793316
    12 = 64.733302669743793
    13 = 97.67954817102779
    14 = 46.725153243037099
}

foreach ($r in $dataRows.Keys) {
    $ws.Range("L$r").Copy()
    $ws.Range("M$r").PasteSpecial($xlPasteFormats)
    $ws.Range("M$r").Value = $dataRows[$r]
}

# --- Refresh the saved cursor position recorded in the sheet view ---
$ws.Range("L19").Select()
